$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("BoM")
$ws.Range("D9").Value = "C22 C26"
$ws.Range("D10").Value = "C29 C34"
$ws.Range("D12").Value = "C1 C11 C13 C15 C16 C17 C18 C19 C20 C21 C23 C25 C30 C35"
$ws.Range("D13").Value = "C12 C14"
$ws.Range("D14").Value = "C2 C3 C4 C5 C6 C7 C8 C9 C10 C24 C27 C28 C32 C36 C37"
$ws.Range("D15").Value = "C33"
$ws.Range("D17").Value = "D7"
$ws.Range("D18").Value = "D6"
$ws.Range("D19").Value = "D2 D3"
$ws.Range("D21").Value = "D4 D5"
$ws.Range("D24").Value = "J27"
$ws.Range("D26").Value = "J5 J8 J18 J19 J20 J22"
$ws.Range("D27").Value = "J9 J10 J13 J15 J17 J21 J23 J24"
$ws.Range("D28").Value = "J14"
$ws.Range("D29").Value = "J28"
$ws.Range("D30").Value = "J11"
$ws.Range("D34").Value = "R5 R13"
$ws.Range("D35").Value = "R6 R10"
$ws.Range("D36").Value = "R4"
$ws.Range("D37").Value = "R3 R7"
$ws.Range("D38").Value = "R1 R2 R8 R11 R16"
$ws.Range("D39").Value = "R15 R17"
$ws.Range("D40").Value = "R12"
$ws.Range("D41").Value = "R14"
$ws.Range("D43").Value = "SW5 SW6 SW7 SW8 SW9 SW10"
$ws.Range("D44").Value = "U1 U7"
$ws.Range("D47").Value = "U3"
$ws.Range("D48").Value = "U2"
$ws.Range("D49").Value = "U4"

$ws = $wb.Worksheets.Item("DNF")
$ws.Range("D10").Value = "J12 J16"
$ws.Range("D11").Value = "J25 J26"
$ws.Range("D13").Value = "R9"
$ws.Range("D14").Value = "SW3 SW4"

$ws = $wb.Worksheets.Item("Costs")
$ws.Range("A10").Value = "C22 C26"
$ws.Range("A11").Value = "C29 C34"
$ws.Range("A13").Value = "C1 C11 C13 C15 C16 C17 C18 C19 C20 C21 C23 C25 C30 C35"
$ws.Range("A14").Value = "C12 C14"
$ws.Range("A15").Value = "C2 C3 C4 C5 C6 C7 C8 C9 C10 C24 C27 C28 C32 C36 C37"
$ws.Range("A16").Value = "C33"
$ws.Range("A18").Value = "D7"
$ws.Range("A19").Value = "D6"
$ws.Range("A20").Value = "D2 D3"
$ws.Range("A22").Value = "D4 D5"
$ws.Range("A25").Value = "J27"
$ws.Range("A27").Value = "J5 J8 J18 J19 J20 J22"
$ws.Range("A28").Value = "J9 J10 J13 J15 J17 J21 J23 J24"
$ws.Range("A29").Value = "J14"
$ws.Range("A30").Value = "J28"
$ws.Range("A31").Value = "J11"
$ws.Range("A35").Value = "R5 R13"
$ws.Range("A36").Value = "R6 R10"
$ws.Range("A37").Value = "R4"
$ws.Range("A38").Value = "R3 R7"
$ws.Range("A39").Value = "R1 R2 R8 R11 R16"
$ws.Range("A40").Value = "R15 R17"
$ws.Range("A41").Value = "R12"
$ws.Range("A42").Value = "R14"
$ws.Range("A44").Value = "SW5 SW6 SW7 SW8 SW9 SW10"
$ws.Range("A45").Value = "U1 U7"
$ws.Range("A48").Value = "U3"
$ws.Range("A49").Value = "U2"
$ws.Range("A50").Value = "U4"
$ws.Range("B54").Value = "2023-12-04 18:11:45"

$ws = $wb.Worksheets.Item("Costs (DNF)")
$ws.Range("A11").Value = "J12 J16"
$ws.Range("A12").Value = "J25 J26"
$ws.Range("A14").Value = "R9"
$ws.Range("A15").Value = "SW3 SW4"
$ws.Range("B18").Value = "2023-12-04 18:11:45"
